$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 is a blank separator row; only H21/I21 carry the date-style formatting (no values).
$ws.Cells.Item(20, 8).Copy()
$ws.Cells.Item(21, 8).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(20, 9).Copy()
$ws.Cells.Item(21, 9).PasteSpecial(-4122)   # xlPasteFormats

# Row 22 is the new data row for LeetCode 594 "Longest Harmonious Subsequence".
$ws.Cells.Item(22, 1).Value = 594
$ws.Cells.Item(22, 2).Value = "Longest Harmonious Subsequence"
$ws.Cells.Item(22, 3).Value = "#array #hash-table #sliding-window #sorting  #counting "
$ws.Cells.Item(22, 4).Value = "easy"
$ws.Cells.Item(22, 5).Value = 1
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 10

$ws.Cells.Item(20, 8).Copy()
$ws.Cells.Item(22, 8).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(22, 8).Value = 45838

$ws.Cells.Item(20, 9).Copy()
$ws.Cells.Item(22, 9).PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(22, 9).Value = 45838

# Match the row height used by the commit for the new data row.
$ws.Rows.Item(22).RowHeight = 68

# Update the view state to match the committed sheet (scrolled down to show the new rows).
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("I22").Select()

